# Update "想去人数" (want-to-go count) figures in column F across the
# "展览" (sheet 1), "演出" (sheet 2) and "全部类型" (sheet 4) worksheets to
# match freshly re-scraped totals (gh-pages data refresh @ 456a3b4).

$wb = $excel.ActiveWorkbook

# ---- Sheet 1: 展览 (Exhibitions) ----
$ws1 = $wb.Worksheets.Item(1)
$ws1.Cells.Item(2, 6).Value = 4649
$ws1.Cells.Item(3, 6).Value = 1848
$ws1.Cells.Item(4, 6).Value = 139
$ws1.Cells.Item(6, 6).Value = 3124
$ws1.Cells.Item(8, 6).Value = 590
$ws1.Cells.Item(9, 6).Value = 270
$ws1.Cells.Item(10, 6).Value = 632
$ws1.Cells.Item(12, 6).Value = 534
$ws1.Cells.Item(13, 6).Value = 383
$ws1.Cells.Item(16, 6).Value = 1346
$ws1.Cells.Item(18, 6).Value = 1618
$ws1.Cells.Item(19, 6).Value = 11
$ws1.Cells.Item(22, 6).Value = 10
$ws1.Cells.Item(23, 6).Value = 45
$ws1.Cells.Item(24, 6).Value = 537
$ws1.Cells.Item(26, 6).Value = 50
$ws1.Cells.Item(27, 6).Value = 104
$ws1.Cells.Item(28, 6).Value = 4
$ws1.Cells.Item(30, 6).Value = 28
$ws1.Cells.Item(32, 6).Value = 3857
$ws1.Cells.Item(33, 6).Value = 3
$ws1.Cells.Item(34, 6).Value = 764
$ws1.Cells.Item(35, 6).Value = 75
$ws1.Cells.Item(36, 6).Value = 893
$ws1.Cells.Item(37, 6).Value = 59
$ws1.Cells.Item(38, 6).Value = 1847

# ---- Sheet 2: 演出 (Performances) ----
$ws2 = $wb.Worksheets.Item(2)
$ws2.Cells.Item(3, 6).Value = 44

# ---- Sheet 4: 全部类型 (All types) ----
$ws4 = $wb.Worksheets.Item(4)
$ws4.Cells.Item(2, 6).Value = 4649
$ws4.Cells.Item(3, 6).Value = 1848
$ws4.Cells.Item(4, 6).Value = 139
$ws4.Cells.Item(6, 6).Value = 3124
$ws4.Cells.Item(8, 6).Value = 590
$ws4.Cells.Item(9, 6).Value = 270
$ws4.Cells.Item(10, 6).Value = 632
$ws4.Cells.Item(12, 6).Value = 534
$ws4.Cells.Item(14, 6).Value = 383
$ws4.Cells.Item(17, 6).Value = 1346
$ws4.Cells.Item(19, 6).Value = 1618
$ws4.Cells.Item(20, 6).Value = 11
$ws4.Cells.Item(23, 6).Value = 10
$ws4.Cells.Item(24, 6).Value = 45
$ws4.Cells.Item(25, 6).Value = 537
$ws4.Cells.Item(27, 6).Value = 50
$ws4.Cells.Item(28, 6).Value = 104
$ws4.Cells.Item(29, 6).Value = 4
$ws4.Cells.Item(31, 6).Value = 28
$ws4.Cells.Item(33, 6).Value = 3858
$ws4.Cells.Item(34, 6).Value = 44
$ws4.Cells.Item(35, 6).Value = 3
$ws4.Cells.Item(36, 6).Value = 764
$ws4.Cells.Item(37, 6).Value = 75
$ws4.Cells.Item(38, 6).Value = 893
$ws4.Cells.Item(39, 6).Value = 59
$ws4.Cells.Item(40, 6).Value = 1847
